$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 800
$ws.Range("A800").Value = 10
$ws.Range("B800").Value = "Vega Modelo de Temuco"
$ws.Range("C800").Value = "La Araucanía"
$ws.Range("D800").Value = 45265
$ws.Range("E800").Value = 9
$ws.Range("F800").Value = "Fruta"
$ws.Range("G800").Value = 100108
$ws.Range("H800").Value = "Tropicales y subtropicales"
$ws.Range("I800").Value = 100108005
$ws.Range("J800").Value = "Piña"
$ws.Range("K800").Value = "Caramelo"
$ws.Range("L800").Value = "Especial"
$ws.Range("M800").Value = 100
$ws.Range("N800").Value = 24000
$ws.Range("O800").Value = 24000
$ws.Range("P800").Value = 24000
$ws.Range("Q800").Value = "`$/caja 10 unidades"
$ws.Range("R800").Value = "Ecuador"
$ws.Range("S800").Value = 2400
$ws.Range("T800").Value = 10

# Row 801
$ws.Range("A801").Value = 10
$ws.Range("B801").Value = "Vega Modelo de Temuco"
$ws.Range("C801").Value = "La Araucanía"
$ws.Range("D801").Value = 45265
$ws.Range("E801").Value = 9
$ws.Range("F801").Value = "Fruta"
$ws.Range("G801").Value = 100108
$ws.Range("H801").Value = "Tropicales y subtropicales"
$ws.Range("I801").Value = 100108005
$ws.Range("J801").Value = "Piña"
$ws.Range("K801").Value = "Caramelo"
$ws.Range("L801").Value = "Primera"
$ws.Range("M801").Value = 40
$ws.Range("N801").Value = 24000
$ws.Range("O801").Value = 24000
$ws.Range("P801").Value = 24000
$ws.Range("Q801").Value = "`$/caja 12 unidades"
$ws.Range("R801").Value = "Ecuador"
$ws.Range("S801").Value = 2000
$ws.Range("T801").Value = 12

# Row 802
$ws.Range("A802").Value = 10
$ws.Range("B802").Value = "Vega Modelo de Temuco"
$ws.Range("C802").Value = "La Araucanía"
$ws.Range("D802").Value = 45222
$ws.Range("E802").Value = 9
$ws.Range("F802").Value = "Fruta"
$ws.Range("G802").Value = 100108
$ws.Range("H802").Value = "Tropicales y subtropicales"
$ws.Range("I802").Value = 100108005
$ws.Range("J802").Value = "Piña"
$ws.Range("K802").Value = "Caramelo"
$ws.Range("L802").Value = "Primera"
$ws.Range("M802").Value = 85
$ws.Range("N802").Value = 24000
$ws.Range("O802").Value = 24000
$ws.Range("P802").Value = 24000
$ws.Range("Q802").Value = "`$/caja 12 unidades"
$ws.Range("R802").Value = "Ecuador"
$ws.Range("S802").Value = 2000
$ws.Range("T802").Value = 12

# Row 803
$ws.Range("A803").Value = 10
$ws.Range("B803").Value = "Vega Modelo de Temuco"
$ws.Range("C803").Value = "La Araucanía"
$ws.Range("D803").Value = 44592
$ws.Range("E803").Value = 9
$ws.Range("F803").Value = "Fruta"
$ws.Range("G803").Value = 100108
$ws.Range("H803").Value = "Tropicales y subtropicales"
$ws.Range("I803").Value = 100108005
$ws.Range("J803").Value = "Piña"
$ws.Range("K803").Value = "Caramelo"
$ws.Range("L803").Value = "Primera"
$ws.Range("M803").Value = 65
$ws.Range("N803").Value = 19000
$ws.Range("O803").Value = 19000
$ws.Range("P803").Value = 19000
$ws.Range("Q803").Value = "`$/caja 12 unidades"
$ws.Range("R803").Value = "Ecuador"
$ws.Range("S803").Value = 1583
$ws.Range("T803").Value = 12

# Row 804
$ws.Range("A804").Value = 10
$ws.Range("B804").Value = "Vega Modelo de Temuco"
$ws.Range("C804").Value = "La Araucanía"
$ws.Range("D804").Value = 45128
$ws.Range("E804").Value = 9
$ws.Range("F804").Value = "Fruta"
$ws.Range("G804").Value = 100108
$ws.Range("H804").Value = "Tropicales y subtropicales"
$ws.Range("I804").Value = 100108005
$ws.Range("J804").Value = "Piña"
$ws.Range("K804").Value = "Caramelo"
$ws.Range("L804").Value = "Segunda"
$ws.Range("M804").Value = 65
$ws.Range("N804").Value = 28000
$ws.Range("O804").Value = 28000
$ws.Range("P804").Value = 28000
$ws.Range("Q804").Value = "`$/caja 14 unidades"
$ws.Range("R804").Value = "Ecuador"
$ws.Range("S804").Value = 2000
$ws.Range("T804").Value = 14

# Row 805
$ws.Range("A805").Value = 10
$ws.Range("B805").Value = "Vega Modelo de Temuco"
$ws.Range("C805").Value = "La Araucanía"
$ws.Range("D805").Value = 44714
$ws.Range("E805").Value = 9
$ws.Range("F805").Value = "Fruta"
$ws.Range("G805").Value = 100108
$ws.Range("H805").Value = "Tropicales y subtropicales"
$ws.Range("I805").Value = 100108005
$ws.Range("J805").Value = "Piña"
$ws.Range("K805").Value = "Caramelo"
$ws.Range("L805").Value = "Primera"
$ws.Range("M805").Value = 200
$ws.Range("N805").Value = 20000
$ws.Range("O805").Value = 20000
$ws.Range("P805").Value = 20000
$ws.Range("Q805").Value = "`$/caja 12 unidades"
$ws.Range("R805").Value = "Ecuador"
$ws.Range("S805").Value = 1667
$ws.Range("T805").Value = 12

# Row 806
$ws.Range("A806").Value = 10
$ws.Range("B806").Value = "Vega Modelo de Temuco"
$ws.Range("C806").Value = "La Araucanía"
$ws.Range("D806").Value = 44714
$ws.Range("E806").Value = 9
$ws.Range("F806").Value = "Fruta"
$ws.Range("G806").Value = 100108
$ws.Range("H806").Value = "Tropicales y subtropicales"
$ws.Range("I806").Value = 100108005
$ws.Range("J806").Value = "Piña"
$ws.Range("K806").Value = "Caramelo"
$ws.Range("L806").Value = "Segunda"
$ws.Range("M806").Value = 300
$ws.Range("N806").Value = 10000
$ws.Range("O806").Value = 10000
$ws.Range("P806").Value = 10000
$ws.Range("Q806").Value = "`$/caja 7 unidades"
$ws.Range("R806").Value = "Ecuador"
$ws.Range("S806").Value = 1429
$ws.Range("T806").Value = 7

# Row 807
$ws.Range("A807").Value = 10
$ws.Range("B807").Value = "Vega Modelo de Temuco"
$ws.Range("C807").Value = "La Araucanía"
$ws.Range("D807").Value = 44761
$ws.Range("E807").Value = 9
$ws.Range("F807").Value = "Fruta"
$ws.Range("G807").Value = 100108
$ws.Range("H807").Value = "Tropicales y subtropicales"
$ws.Range("I807").Value = 100108005
$ws.Range("J807").Value = "Piña"
$ws.Range("K807").Value = "Caramelo"
$ws.Range("L807").Value = "Primera"
$ws.Range("M807").Value = 100
$ws.Range("N807").Value = 23000
$ws.Range("O807").Value = 23000
$ws.Range("P807").Value = 23000
$ws.Range("Q807").Value = "`$/caja 12 unidades"
$ws.Range("R807").Value = "Ecuador"
$ws.Range("S807").Value = 1917
$ws.Range("T807").Value = 12

# Row 808
$ws.Range("A808").Value = 10
$ws.Range("B808").Value = "Vega Modelo de Temuco"
$ws.Range("C808").Value = "La Araucanía"
$ws.Range("D808").Value = 44363
$ws.Range("E808").Value = 9
$ws.Range("F808").Value = "Fruta"
$ws.Range("G808").Value = 100108
$ws.Range("H808").Value = "Tropicales y subtropicales"
$ws.Range("I808").Value = 100108005
$ws.Range("J808").Value = "Piña"
$ws.Range("K808").Value = "Caramelo"
$ws.Range("L808").Value = "Primera"
$ws.Range("M808").Value = 80
$ws.Range("N808").Value = 17000
$ws.Range("O808").Value = 17000
$ws.Range("P808").Value = 17000
$ws.Range("Q808").Value = "`$/caja 12 unidades"
$ws.Range("R808").Value = "Ecuador"
$ws.Range("S808").Value = 1417
$ws.Range("T808").Value = 12

# Row 809
$ws.Range("A809").Value = 10
$ws.Range("B809").Value = "Vega Modelo de Temuco"
$ws.Range("C809").Value = "La Araucanía"
$ws.Range("D809").Value = 44222
$ws.Range("E809").Value = 9
$ws.Range("F809").Value = "Fruta"
$ws.Range("G809").Value = 100108
$ws.Range("H809").Value = "Tropicales y subtropicales"
$ws.Range("I809").Value = 100108005
$ws.Range("J809").Value = "Piña"
$ws.Range("K809").Value = "Caramelo"
$ws.Range("L809").Value = "Primera"
$ws.Range("M809").Value = 170
$ws.Range("N809").Value = 16000
$ws.Range("O809").Value = 17000
$ws.Range("P809").Value = 16471
$ws.Range("Q809").Value = "`$/caja 12 unidades"
$ws.Range("R809").Value = "Ecuador"
$ws.Range("S809").Value = 1373
$ws.Range("T809").Value = 12

# Row 810
$ws.Range("A810").Value = 10
$ws.Range("B810").Value = "Vega Modelo de Temuco"
$ws.Range("C810").Value = "La Araucanía"
$ws.Range("D810").Value = 44222
$ws.Range("E810").Value = 9
$ws.Range("F810").Value = "Fruta"
$ws.Range("G810").Value = 100108
$ws.Range("H810").Value = "Tropicales y subtropicales"
$ws.Range("I810").Value = 100108005
$ws.Range("J810").Value = "Piña"
$ws.Range("K810").Value = "Caramelo"
$ws.Range("L810").Value = "Segunda"
$ws.Range("M810").Value = 110
$ws.Range("N810").Value = 15000
$ws.Range("O810").Value = 16000
$ws.Range("P810").Value = 15455
$ws.Range("Q810").Value = "`$/caja 14 unidades"
$ws.Range("R810").Value = "Ecuador"
$ws.Range("S810").Value = 1104
$ws.Range("T810").Value = 14

# Row 811
$ws.Range("A811").Value = 10
$ws.Range("B811").Value = "Vega Modelo de Temuco"
$ws.Range("C811").Value = "La Araucanía"
$ws.Range("D811").Value = 44271
$ws.Range("E811").Value = 9
$ws.Range("F811").Value = "Fruta"
$ws.Range("G811").Value = 100108
$ws.Range("H811").Value = "Tropicales y subtropicales"
$ws.Range("I811").Value = 100108005
$ws.Range("J811").Value = "Piña"
$ws.Range("K811").Value = "Caramelo"
$ws.Range("L811").Value = "Segunda"
$ws.Range("M811").Value = 75
$ws.Range("N811").Value = 15000
$ws.Range("O811").Value = 16000
$ws.Range("P811").Value = 15467
$ws.Range("Q811").Value = "`$/caja 14 unidades"
$ws.Range("R811").Value = "Ecuador"
$ws.Range("S811").Value = 1105
$ws.Range("T811").Value = 14

# Row 812
$ws.Range("A812").Value = 10
$ws.Range("B812").Value = "Vega Modelo de Temuco"
$ws.Range("C812").Value = "La Araucanía"
$ws.Range("D812").Value = 44880
$ws.Range("E812").Value = 9
$ws.Range("F812").Value = "Fruta"
$ws.Range("G812").Value = 100108
$ws.Range("H812").Value = "Tropicales y subtropicales"
$ws.Range("I812").Value = 100108005
$ws.Range("J812").Value = "Piña"
$ws.Range("K812").Value = "Caramelo"
$ws.Range("L812").Value = "Segunda"
$ws.Range("M812").Value = 65
$ws.Range("N812").Value = 28000
$ws.Range("O812").Value = 28000
$ws.Range("P812").Value = 28000
$ws.Range("Q812").Value = "`$/caja 14 unidades"
$ws.Range("R812").Value = "Ecuador"
$ws.Range("S812").Value = 2000
$ws.Range("T812").Value = 14

# Row 813
$ws.Range("A813").Value = 10
$ws.Range("B813").Value = "Vega Modelo de Temuco"
$ws.Range("C813").Value = "La Araucanía"
$ws.Range("D813").Value = 44650
$ws.Range("E813").Value = 9
$ws.Range("F813").Value = "Fruta"
$ws.Range("G813").Value = 100108
$ws.Range("H813").Value = "Tropicales y subtropicales"
$ws.Range("I813").Value = 100108005
$ws.Range("J813").Value = "Piña"
$ws.Range("K813").Value = "Caramelo"
$ws.Range("L813").Value = "Primera"
$ws.Range("M813").Value = 115
$ws.Range("N813").Value = 17000
$ws.Range("O813").Value = 18000
$ws.Range("P813").Value = 17565
$ws.Range("Q813").Value = "`$/caja 12 unidades"
$ws.Range("R813").Value = "Ecuador"
$ws.Range("S813").Value = 1464
$ws.Range("T813").Value = 12

# Row 814
$ws.Range("A814").Value = 10
$ws.Range("B814").Value = "Vega Modelo de Temuco"
$ws.Range("C814").Value = "La Araucanía"
$ws.Range("D814").Value = 44399
$ws.Range("E814").Value = 9
$ws.Range("F814").Value = "Fruta"
$ws.Range("G814").Value = 100108
$ws.Range("H814").Value = "Tropicales y subtropicales"
$ws.Range("I814").Value = 100108005
$ws.Range("J814").Value = "Piña"
$ws.Range("K814").Value = "Caramelo"
$ws.Range("L814").Value = "Primera"
$ws.Range("M814").Value = 55
$ws.Range("N814").Value = 20000
$ws.Range("O814").Value = 20000
$ws.Range("P814").Value = 20000
$ws.Range("Q814").Value = "`$/caja 12 unidades"
$ws.Range("R814").Value = "Ecuador"
$ws.Range("S814").Value = 1667
$ws.Range("T814").Value = 12

# Row 815
$ws.Range("A815").Value = 10
$ws.Range("B815").Value = "Vega Modelo de Temuco"
$ws.Range("C815").Value = "La Araucanía"
$ws.Range("D815").Value = 44399
$ws.Range("E815").Value = 9
$ws.Range("F815").Value = "Fruta"
$ws.Range("G815").Value = 100108
$ws.Range("H815").Value = "Tropicales y subtropicales"
$ws.Range("I815").Value = 100108005
$ws.Range("J815").Value = "Piña"
$ws.Range("K815").Value = "Caramelo"
$ws.Range("L815").Value = "Segunda"
$ws.Range("M815").Value = 75
$ws.Range("N815").Value = 20000
$ws.Range("O815").Value = 20000
$ws.Range("P815").Value = 20000
$ws.Range("Q815").Value = "`$/caja 14 unidades"
$ws.Range("R815").Value = "Ecuador"
$ws.Range("S815").Value = 1429
$ws.Range("T815").Value = 14

# Row 816
$ws.Range("A816").Value = 10
$ws.Range("B816").Value = "Vega Modelo de Temuco"
$ws.Range("C816").Value = "La Araucanía"
$ws.Range("D816").Value = 44161
$ws.Range("E816").Value = 9
$ws.Range("F816").Value = "Fruta"
$ws.Range("G816").Value = 100108
$ws.Range("H816").Value = "Tropicales y subtropicales"
$ws.Range("I816").Value = 100108005
$ws.Range("J816").Value = "Piña"
$ws.Range("K816").Value = "Caramelo"
$ws.Range("L816").Value = "Segunda"
$ws.Range("M816").Value = 90
$ws.Range("N816").Value = 23000
$ws.Range("O816").Value = 25000
$ws.Range("P816").Value = 24222
$ws.Range("Q816").Value = "`$/caja 14 unidades"
$ws.Range("R816").Value = "Ecuador"
$ws.Range("S816").Value = 1730
$ws.Range("T816").Value = 14

# Row 817
$ws.Range("A817").Value = 10
$ws.Range("B817").Value = "Vega Modelo de Temuco"
$ws.Range("C817").Value = "La Araucanía"
$ws.Range("D817").Value = 44161
$ws.Range("E817").Value = 9
$ws.Range("F817").Value = "Fruta"
$ws.Range("G817").Value = 100108
$ws.Range("H817").Value = "Tropicales y subtropicales"
$ws.Range("I817").Value = 100108005
$ws.Range("J817").Value = "Piña"
$ws.Range("K817").Value = "Caramelo"
$ws.Range("L817").Value = "Segunda"
$ws.Range("M817").Value = 85
$ws.Range("N817").Value = 15000
$ws.Range("O817").Value = 15000
$ws.Range("P817").Value = 15000
$ws.Range("Q817").Value = "`$/caja 7 unidades"
$ws.Range("R817").Value = "Ecuador"
$ws.Range("S817").Value = 2143
$ws.Range("T817").Value = 7

# Row 818
$ws.Range("A818").Value = 10
$ws.Range("B818").Value = "Vega Modelo de Temuco"
$ws.Range("C818").Value = "La Araucanía"
$ws.Range("D818").Value = 44435
$ws.Range("E818").Value = 9
$ws.Range("F818").Value = "Fruta"
$ws.Range("G818").Value = 100108
$ws.Range("H818").Value = "Tropicales y subtropicales"
$ws.Range("I818").Value = 100108005
$ws.Range("J818").Value = "Piña"
$ws.Range("K818").Value = "Caramelo"
$ws.Range("L818").Value = "Primera"
$ws.Range("M818").Value = 590
$ws.Range("N818").Value = 19000
$ws.Range("O818").Value = 21000
$ws.Range("P818").Value = 19780
$ws.Range("Q818").Value = "`$/caja 12 unidades"
$ws.Range("R818").Value = "Ecuador"
$ws.Range("S818").Value = 1648
$ws.Range("T818").Value = 12

# Row 819
$ws.Range("A819").Value = 10
$ws.Range("B819").Value = "Vega Modelo de Temuco"
$ws.Range("C819").Value = "La Araucanía"
$ws.Range("D819").Value = 44435
$ws.Range("E819").Value = 9
$ws.Range("F819").Value = "Fruta"
$ws.Range("G819").Value = 100108
$ws.Range("H819").Value = "Tropicales y subtropicales"
$ws.Range("I819").Value = 100108005
$ws.Range("J819").Value = "Piña"
$ws.Range("K819").Value = "Caramelo"
$ws.Range("L819").Value = "Segunda"
$ws.Range("M819").Value = 180
$ws.Range("N819").Value = 17000
$ws.Range("O819").Value = 21000
$ws.Range("P819").Value = 19444
$ws.Range("Q819").Value = "`$/caja 14 unidades"
$ws.Range("R819").Value = "Ecuador"
$ws.Range("S819").Value = 1389
$ws.Range("T819").Value = 14

# Row 820
$ws.Range("A820").Value = 10
$ws.Range("B820").Value = "Vega Modelo de Temuco"
$ws.Range("C820").Value = "La Araucanía"
$ws.Range("D820").Value = 44609
$ws.Range("E820").Value = 9
$ws.Range("F820").Value = "Fruta"
$ws.Range("G820").Value = 100108
$ws.Range("H820").Value = "Tropicales y subtropicales"
$ws.Range("I820").Value = 100108005
$ws.Range("J820").Value = "Piña"
$ws.Range("K820").Value = "Caramelo"
$ws.Range("L820").Value = "Primera"
$ws.Range("M820").Value = 400
$ws.Range("N820").Value = 19000
$ws.Range("O820").Value = 20000
$ws.Range("P820").Value = 19500
$ws.Range("Q820").Value = "`$/caja 12 unidades"
$ws.Range("R820").Value = "Ecuador"
$ws.Range("S820").Value = 1625
$ws.Range("T820").Value = 12

# Row 821
$ws.Range("A821").Value = 10
$ws.Range("B821").Value = "Vega Modelo de Temuco"
$ws.Range("C821").Value = "La Araucanía"
$ws.Range("D821").Value = 44609
$ws.Range("E821").Value = 9
$ws.Range("F821").Value = "Fruta"
$ws.Range("G821").Value = 100108
$ws.Range("H821").Value = "Tropicales y subtropicales"
$ws.Range("I821").Value = 100108005
$ws.Range("J821").Value = "Piña"
$ws.Range("K821").Value = "Caramelo"
$ws.Range("L821").Value = "Segunda"
$ws.Range("M821").Value = 100
$ws.Range("N821").Value = 17000
$ws.Range("O821").Value = 17000
$ws.Range("P821").Value = 17000
$ws.Range("Q821").Value = "`$/caja 14 unidades"
$ws.Range("R821").Value = "Ecuador"
$ws.Range("S821").Value = 1214
$ws.Range("T821").Value = 14

# Row 822
$ws.Range("A822").Value = 10
$ws.Range("B822").Value = "Vega Modelo de Temuco"
$ws.Range("C822").Value = "La Araucanía"
$ws.Range("D822").Value = 44509
$ws.Range("E822").Value = 9
$ws.Range("F822").Value = "Fruta"
$ws.Range("G822").Value = 100108
$ws.Range("H822").Value = "Tropicales y subtropicales"
$ws.Range("I822").Value = 100108005
$ws.Range("J822").Value = "Piña"
$ws.Range("K822").Value = "Caramelo"
$ws.Range("L822").Value = "Primera"
$ws.Range("M822").Value = 55
$ws.Range("N822").Value = 21000
$ws.Range("O822").Value = 21000
$ws.Range("P822").Value = 21000
$ws.Range("Q822").Value = "`$/caja 12 unidades"
$ws.Range("R822").Value = "Ecuador"
$ws.Range("S822").Value = 1750
$ws.Range("T822").Value = 12

# Row 823
$ws.Range("A823").Value = 10
$ws.Range("B823").Value = "Vega Modelo de Temuco"
$ws.Range("C823").Value = "La Araucanía"
$ws.Range("D823").Value = 44509
$ws.Range("E823").Value = 9
$ws.Range("F823").Value = "Fruta"
$ws.Range("G823").Value = 100108
$ws.Range("H823").Value = "Tropicales y subtropicales"
$ws.Range("I823").Value = 100108005
$ws.Range("J823").Value = "Piña"
$ws.Range("K823").Value = "Caramelo"
$ws.Range("L823").Value = "Segunda"
$ws.Range("M823").Value = 45
$ws.Range("N823").Value = 20000
$ws.Range("O823").Value = 20000
$ws.Range("P823").Value = 20000
$ws.Range("Q823").Value = "`$/caja 14 unidades"
$ws.Range("R823").Value = "Ecuador"
$ws.Range("S823").Value = 1429
$ws.Range("T823").Value = 14

# Row 824
$ws.Range("A824").Value = 10
$ws.Range("B824").Value = "Vega Modelo de Temuco"
$ws.Range("C824").Value = "La Araucanía"
$ws.Range("D824").Value = 45050
$ws.Range("E824").Value = 9
$ws.Range("F824").Value = "Fruta"
$ws.Range("G824").Value = 100108
$ws.Range("H824").Value = "Tropicales y subtropicales"
$ws.Range("I824").Value = 100108005
$ws.Range("J824").Value = "Piña"
$ws.Range("K824").Value = "Caramelo"
$ws.Range("L824").Value = "Primera"
$ws.Range("M824").Value = 250
$ws.Range("N824").Value = 20000
$ws.Range("O824").Value = 20000
$ws.Range("P824").Value = 20000
$ws.Range("Q824").Value = "`$/caja 12 unidades"
$ws.Range("R824").Value = "Ecuador"
$ws.Range("S824").Value = 1667
$ws.Range("T824").Value = 12

# Row 825
$ws.Range("A825").Value = 10
$ws.Range("B825").Value = "Vega Modelo de Temuco"
$ws.Range("C825").Value = "La Araucanía"
$ws.Range("D825").Value = 45180
$ws.Range("E825").Value = 9
$ws.Range("F825").Value = "Fruta"
$ws.Range("G825").Value = 100108
$ws.Range("H825").Value = "Tropicales y subtropicales"
$ws.Range("I825").Value = 100108005
$ws.Range("J825").Value = "Piña"
$ws.Range("K825").Value = "Caramelo"
$ws.Range("L825").Value = "Primera"
$ws.Range("M825").Value = 150
$ws.Range("N825").Value = 25000
$ws.Range("O825").Value = 25000
$ws.Range("P825").Value = 25000
$ws.Range("Q825").Value = "`$/caja 12 unidades"
$ws.Range("R825").Value = "Ecuador"
$ws.Range("S825").Value = 2083
$ws.Range("T825").Value = 12

# Row 826
$ws.Range("A826").Value = 10
$ws.Range("B826").Value = "Vega Modelo de Temuco"
$ws.Range("C826").Value = "La Araucanía"
$ws.Range("D826").Value = 44824
$ws.Range("E826").Value = 9
$ws.Range("F826").Value = "Fruta"
$ws.Range("G826").Value = 100108
$ws.Range("H826").Value = "Tropicales y subtropicales"
$ws.Range("I826").Value = 100108005
$ws.Range("J826").Value = "Piña"
$ws.Range("K826").Value = "Caramelo"
$ws.Range("L826").Value = "Segunda"
$ws.Range("M826").Value = 45
$ws.Range("N826").Value = 22000
$ws.Range("O826").Value = 23000
$ws.Range("P826").Value = 22444
$ws.Range("Q826").Value = "`$/caja 14 unidades"
$ws.Range("R826").Value = "Ecuador"
$ws.Range("S826").Value = 1603
$ws.Range("T826").Value = 14

# Row 827
$ws.Range("A827").Value = 10
$ws.Range("B827").Value = "Vega Modelo de Temuco"
$ws.Range("C827").Value = "La Araucanía"
$ws.Range("D827").Value = 44428
$ws.Range("E827").Value = 9
$ws.Range("F827").Value = "Fruta"
$ws.Range("G827").Value = 100108
$ws.Range("H827").Value = "Tropicales y subtropicales"
$ws.Range("I827").Value = 100108005
$ws.Range("J827").Value = "Piña"
$ws.Range("K827").Value = "Caramelo"
$ws.Range("L827").Value = "Primera"
$ws.Range("M827").Value = 65
$ws.Range("N827").Value = 19000
$ws.Range("O827").Value = 20000
$ws.Range("P827").Value = 19538
$ws.Range("Q827").Value = "`$/caja 12 unidades"
$ws.Range("R827").Value = "Ecuador"
$ws.Range("S827").Value = 1628
$ws.Range("T827").Value = 12

# Row 828
$ws.Range("A828").Value = 10
$ws.Range("B828").Value = "Vega Modelo de Temuco"
$ws.Range("C828").Value = "La Araucanía"
$ws.Range("D828").Value = 44428
$ws.Range("E828").Value = 9
$ws.Range("F828").Value = "Fruta"
$ws.Range("G828").Value = 100108
$ws.Range("H828").Value = "Tropicales y subtropicales"
$ws.Range("I828").Value = 100108005
$ws.Range("J828").Value = "Piña"
$ws.Range("K828").Value = "Caramelo"
$ws.Range("L828").Value = "Segunda"
$ws.Range("M828").Value = 65
$ws.Range("N828").Value = 19000
$ws.Range("O828").Value = 19000
$ws.Range("P828").Value = 19000
$ws.Range("Q828").Value = "`$/caja 14 unidades"
$ws.Range("R828").Value = "Ecuador"
$ws.Range("S828").Value = 1357
$ws.Range("T828").Value = 14

# Row 829
$ws.Range("A829").Value = 10
$ws.Range("B829").Value = "Vega Modelo de Temuco"
$ws.Range("C829").Value = "La Araucanía"
$ws.Range("D829").Value = 44547
$ws.Range("E829").Value = 9
$ws.Range("F829").Value = "Fruta"
$ws.Range("G829").Value = 100108
$ws.Range("H829").Value = "Tropicales y subtropicales"
$ws.Range("I829").Value = 100108005
$ws.Range("J829").Value = "Piña"
$ws.Range("K829").Value = "Caramelo"
$ws.Range("L829").Value = "Primera"
$ws.Range("M829").Value = 115
$ws.Range("N829").Value = 20000
$ws.Range("O829").Value = 20000
$ws.Range("P829").Value = 20000
$ws.Range("Q829").Value = "`$/caja 12 unidades"
$ws.Range("R829").Value = "Ecuador"
$ws.Range("S829").Value = 1667
$ws.Range("T829").Value = 12

# Row 830
$ws.Range("A830").Value = 10
$ws.Range("B830").Value = "Vega Modelo de Temuco"
$ws.Range("C830").Value = "La Araucanía"
$ws.Range("D830").Value = 44547
$ws.Range("E830").Value = 9
$ws.Range("F830").Value = "Fruta"
$ws.Range("G830").Value = 100108
$ws.Range("H830").Value = "Tropicales y subtropicales"
$ws.Range("I830").Value = 100108005
$ws.Range("J830").Value = "Piña"
$ws.Range("K830").Value = "Caramelo"
$ws.Range("L830").Value = "Segunda"
$ws.Range("M830").Value = 155
$ws.Range("N830").Value = 20000
$ws.Range("O830").Value = 20000
$ws.Range("P830").Value = 20000
$ws.Range("Q830").Value = "`$/caja 14 unidades"
$ws.Range("R830").Value = "Ecuador"
$ws.Range("S830").Value = 1429
$ws.Range("T830").Value = 14

# Row 831
$ws.Range("A831").Value = 10
$ws.Range("B831").Value = "Vega Modelo de Temuco"
$ws.Range("C831").Value = "La Araucanía"
$ws.Range("D831").Value = 44539
$ws.Range("E831").Value = 9
$ws.Range("F831").Value = "Fruta"
$ws.Range("G831").Value = 100108
$ws.Range("H831").Value = "Tropicales y subtropicales"
$ws.Range("I831").Value = 100108005
$ws.Range("J831").Value = "Piña"
$ws.Range("K831").Value = "Caramelo"
$ws.Range("L831").Value = "Primera"
$ws.Range("M831").Value = 65
$ws.Range("N831").Value = 21000
$ws.Range("O831").Value = 21000
$ws.Range("P831").Value = 21000
$ws.Range("Q831").Value = "`$/caja 14 unidades"
$ws.Range("R831").Value = "Ecuador"
$ws.Range("S831").Value = 1500
$ws.Range("T831").Value = 14

# Row 832
$ws.Range("A832").Value = 10
$ws.Range("B832").Value = "Vega Modelo de Temuco"
$ws.Range("C832").Value = "La Araucanía"
$ws.Range("D832").Value = 44894
$ws.Range("E832").Value = 9
$ws.Range("F832").Value = "Fruta"
$ws.Range("G832").Value = 100108
$ws.Range("H832").Value = "Tropicales y subtropicales"
$ws.Range("I832").Value = 100108005
$ws.Range("J832").Value = "Piña"
$ws.Range("K832").Value = "Caramelo"
$ws.Range("L832").Value = "Segunda"
$ws.Range("M832").Value = 35
$ws.Range("N832").Value = 32000
$ws.Range("O832").Value = 32000
$ws.Range("P832").Value = 32000
$ws.Range("Q832").Value = "`$/caja 14 unidades"
$ws.Range("R832").Value = "Ecuador"
$ws.Range("S832").Value = 2286
$ws.Range("T832").Value = 14

# Row 833
$ws.Range("A833").Value = 10
$ws.Range("B833").Value = "Vega Modelo de Temuco"
$ws.Range("C833").Value = "La Araucanía"
$ws.Range("D833").Value = 44193
$ws.Range("E833").Value = 9
$ws.Range("F833").Value = "Fruta"
$ws.Range("G833").Value = 100108
$ws.Range("H833").Value = "Tropicales y subtropicales"
$ws.Range("I833").Value = 100108005
$ws.Range("J833").Value = "Piña"
$ws.Range("K833").Value = "Caramelo"
$ws.Range("L833").Value = "Segunda"
$ws.Range("M833").Value = 100
$ws.Range("N833").Value = 17000
$ws.Range("O833").Value = 18000
$ws.Range("P833").Value = 17550
$ws.Range("Q833").Value = "`$/caja 14 unidades"
$ws.Range("R833").Value = "Ecuador"
$ws.Range("S833").Value = 1254
$ws.Range("T833").Value = 14

# Row 834
$ws.Range("A834").Value = 10
$ws.Range("B834").Value = "Vega Modelo de Temuco"
$ws.Range("C834").Value = "La Araucanía"
$ws.Range("D834").Value = 44917
$ws.Range("E834").Value = 9
$ws.Range("F834").Value = "Fruta"
$ws.Range("G834").Value = 100108
$ws.Range("H834").Value = "Tropicales y subtropicales"
$ws.Range("I834").Value = 100108005
$ws.Range("J834").Value = "Piña"
$ws.Range("K834").Value = "Caramelo"
$ws.Range("L834").Value = "Primera"
$ws.Range("M834").Value = 250
$ws.Range("N834").Value = 18000
$ws.Range("O834").Value = 20000
$ws.Range("P834").Value = 18800
$ws.Range("Q834").Value = "`$/caja 12 unidades"
$ws.Range("R834").Value = "Ecuador"
$ws.Range("S834").Value = 1567
$ws.Range("T834").Value = 12

# Row 835
$ws.Range("A835").Value = 10
$ws.Range("B835").Value = "Vega Modelo de Temuco"
$ws.Range("C835").Value = "La Araucanía"
$ws.Range("D835").Value = 44631
$ws.Range("E835").Value = 9
$ws.Range("F835").Value = "Fruta"
$ws.Range("G835").Value = 100108
$ws.Range("H835").Value = "Tropicales y subtropicales"
$ws.Range("I835").Value = 100108005
$ws.Range("J835").Value = "Piña"
$ws.Range("K835").Value = "Caramelo"
$ws.Range("L835").Value = "Primera"
$ws.Range("M835").Value = 80
$ws.Range("N835").Value = 18000
$ws.Range("O835").Value = 19000
$ws.Range("P835").Value = 18500
$ws.Range("Q835").Value = "`$/caja 12 unidades"
$ws.Range("R835").Value = "Ecuador"
$ws.Range("S835").Value = 1542
$ws.Range("T835").Value = 12

# Row 836
$ws.Range("A836").Value = 10
$ws.Range("B836").Value = "Vega Modelo de Temuco"
$ws.Range("C836").Value = "La Araucanía"
$ws.Range("D836").Value = 44631
$ws.Range("E836").Value = 9
$ws.Range("F836").Value = "Fruta"
$ws.Range("G836").Value = 100108
$ws.Range("H836").Value = "Tropicales y subtropicales"
$ws.Range("I836").Value = 100108005
$ws.Range("J836").Value = "Piña"
$ws.Range("K836").Value = "Caramelo"
$ws.Range("L836").Value = "Segunda"
$ws.Range("M836").Value = 60
$ws.Range("N836").Value = 18000
$ws.Range("O836").Value = 19000
$ws.Range("P836").Value = 18500
$ws.Range("Q836").Value = "`$/caja 14 unidades"
$ws.Range("R836").Value = "Ecuador"
$ws.Range("S836").Value = 1321
$ws.Range("T836").Value = 14

# Row 837
$ws.Range("A837").Value = 10
$ws.Range("B837").Value = "Vega Modelo de Temuco"
$ws.Range("C837").Value = "La Araucanía"
$ws.Range("D837").Value = 44727
$ws.Range("E837").Value = 9
$ws.Range("F837").Value = "Fruta"
$ws.Range("G837").Value = 100108
$ws.Range("H837").Value = "Tropicales y subtropicales"
$ws.Range("I837").Value = 100108005
$ws.Range("J837").Value = "Piña"
$ws.Range("K837").Value = "Caramelo"
$ws.Range("L837").Value = "Primera"
$ws.Range("M837").Value = 90
$ws.Range("N837").Value = 20000
$ws.Range("O837").Value = 20000
$ws.Range("P837").Value = 20000
$ws.Range("Q837").Value = "`$/caja 12 unidades"
$ws.Range("R837").Value = "Ecuador"
$ws.Range("S837").Value = 1667
$ws.Range("T837").Value = 12

# Row 838
$ws.Range("A838").Value = 10
$ws.Range("B838").Value = "Vega Modelo de Temuco"
$ws.Range("C838").Value = "La Araucanía"
$ws.Range("D838").Value = 44914
$ws.Range("E838").Value = 9
$ws.Range("F838").Value = "Fruta"
$ws.Range("G838").Value = 100108
$ws.Range("H838").Value = "Tropicales y subtropicales"
$ws.Range("I838").Value = 100108005
$ws.Range("J838").Value = "Piña"
$ws.Range("K838").Value = "Caramelo"
$ws.Range("L838").Value = "Primera"
$ws.Range("M838").Value = 55
$ws.Range("N838").Value = 22000
$ws.Range("O838").Value = 22000
$ws.Range("P838").Value = 22000
$ws.Range("Q838").Value = "`$/caja 12 unidades"
$ws.Range("R838").Value = "Ecuador"
$ws.Range("S838").Value = 1833
$ws.Range("T838").Value = 12

# Row 839
$ws.Range("A839").Value = 10
$ws.Range("B839").Value = "Vega Modelo de Temuco"
$ws.Range("C839").Value = "La Araucanía"
$ws.Range("D839").Value = 44914
$ws.Range("E839").Value = 9
$ws.Range("F839").Value = "Fruta"
$ws.Range("G839").Value = 100108
$ws.Range("H839").Value = "Tropicales y subtropicales"
$ws.Range("I839").Value = 100108005
$ws.Range("J839").Value = "Piña"
$ws.Range("K839").Value = "Caramelo"
$ws.Range("L839").Value = "Segunda"
$ws.Range("M839").Value = 65
$ws.Range("N839").Value = 22000
$ws.Range("O839").Value = 22000
$ws.Range("P839").Value = 22000
$ws.Range("Q839").Value = "`$/caja 14 unidades"
$ws.Range("R839").Value = "Ecuador"
$ws.Range("S839").Value = 1571
$ws.Range("T839").Value = 14

# Row 840
$ws.Range("A840").Value = 10
$ws.Range("B840").Value = "Vega Modelo de Temuco"
$ws.Range("C840").Value = "La Araucanía"
$ws.Range("D840").Value = 44239
$ws.Range("E840").Value = 9
$ws.Range("F840").Value = "Fruta"
$ws.Range("G840").Value = 100108
$ws.Range("H840").Value = "Tropicales y subtropicales"
$ws.Range("I840").Value = 100108005
$ws.Range("J840").Value = "Piña"
$ws.Range("K840").Value = "Caramelo"
$ws.Range("L840").Value = "Especial"
$ws.Range("M840").Value = 80
$ws.Range("N840").Value = 17000
$ws.Range("O840").Value = 18000
$ws.Range("P840").Value = 17562
$ws.Range("Q840").Value = "`$/caja 10 unidades"
$ws.Range("R840").Value = "Ecuador"
$ws.Range("S840").Value = 1756
$ws.Range("T840").Value = 10

# Row 841
$ws.Range("A841").Value = 10
$ws.Range("B841").Value = "Vega Modelo de Temuco"
$ws.Range("C841").Value = "La Araucanía"
$ws.Range("D841").Value = 44239
$ws.Range("E841").Value = 9
$ws.Range("F841").Value = "Fruta"
$ws.Range("G841").Value = 100108
$ws.Range("H841").Value = "Tropicales y subtropicales"
$ws.Range("I841").Value = 100108005
$ws.Range("J841").Value = "Piña"
$ws.Range("K841").Value = "Caramelo"
$ws.Range("L841").Value = "Primera"
$ws.Range("M841").Value = 55
$ws.Range("N841").Value = 16000
$ws.Range("O841").Value = 16000
$ws.Range("P841").Value = 16000
$ws.Range("Q841").Value = "`$/caja 12 unidades"
$ws.Range("R841").Value = "Ecuador"
$ws.Range("S841").Value = 1333
$ws.Range("T841").Value = 12

# Row 842
$ws.Range("A842").Value = 10
$ws.Range("B842").Value = "Vega Modelo de Temuco"
$ws.Range("C842").Value = "La Araucanía"
$ws.Range("D842").Value = 45168
$ws.Range("E842").Value = 9
$ws.Range("F842").Value = "Fruta"
$ws.Range("G842").Value = 100108
$ws.Range("H842").Value = "Tropicales y subtropicales"
$ws.Range("I842").Value = 100108005
$ws.Range("J842").Value = "Piña"
$ws.Range("K842").Value = "Caramelo"
$ws.Range("L842").Value = "Primera"
$ws.Range("M842").Value = 35
$ws.Range("N842").Value = 25000
$ws.Range("O842").Value = 25000
$ws.Range("P842").Value = 25000
$ws.Range("Q842").Value = "`$/caja 12 unidades"
$ws.Range("R842").Value = "Ecuador"
$ws.Range("S842").Value = 2083
$ws.Range("T842").Value = 12

# Row 843
$ws.Range("A843").Value = 10
$ws.Range("B843").Value = "Vega Modelo de Temuco"
$ws.Range("C843").Value = "La Araucanía"
$ws.Range("D843").Value = 44574
$ws.Range("E843").Value = 9
$ws.Range("F843").Value = "Fruta"
$ws.Range("G843").Value = 100108
$ws.Range("H843").Value = "Tropicales y subtropicales"
$ws.Range("I843").Value = 100108005
$ws.Range("J843").Value = "Piña"
$ws.Range("K843").Value = "Caramelo"
$ws.Range("L843").Value = "Primera"
$ws.Range("M843").Value = 400
$ws.Range("N843").Value = 19000
$ws.Range("O843").Value = 19000
$ws.Range("P843").Value = 19000
$ws.Range("Q843").Value = "`$/caja 12 unidades"
$ws.Range("R843").Value = "Ecuador"
$ws.Range("S843").Value = 1583
$ws.Range("T843").Value = 12

# Row 844
$ws.Range("A844").Value = 10
$ws.Range("B844").Value = "Vega Modelo de Temuco"
$ws.Range("C844").Value = "La Araucanía"
$ws.Range("D844").Value = 44574
$ws.Range("E844").Value = 9
$ws.Range("F844").Value = "Fruta"
$ws.Range("G844").Value = 100108
$ws.Range("H844").Value = "Tropicales y subtropicales"
$ws.Range("I844").Value = 100108005
$ws.Range("J844").Value = "Piña"
$ws.Range("K844").Value = "Caramelo"
$ws.Range("L844").Value = "Segunda"
$ws.Range("M844").Value = 100
$ws.Range("N844").Value = 19000
$ws.Range("O844").Value = 19000
$ws.Range("P844").Value = 19000
$ws.Range("Q844").Value = "`$/caja 14 unidades"
$ws.Range("R844").Value = "Ecuador"
$ws.Range("S844").Value = 1357
$ws.Range("T844").Value = 14

# Row 845
$ws.Range("A845").Value = 10
$ws.Range("B845").Value = "Vega Modelo de Temuco"
$ws.Range("C845").Value = "La Araucanía"
$ws.Range("D845").Value = 44574
$ws.Range("E845").Value = 9
$ws.Range("F845").Value = "Fruta"
$ws.Range("G845").Value = 100108
$ws.Range("H845").Value = "Tropicales y subtropicales"
$ws.Range("I845").Value = 100108005
$ws.Range("J845").Value = "Piña"
$ws.Range("K845").Value = "Caramelo"
$ws.Range("L845").Value = "Tercera"
$ws.Range("M845").Value = 80
$ws.Range("N845").Value = 19000
$ws.Range("O845").Value = 19000
$ws.Range("P845").Value = 19000
$ws.Range("Q845").Value = "`$/caja 16 unidades"
$ws.Range("R845").Value = "Ecuador"
$ws.Range("S845").Value = 1188
$ws.Range("T845").Value = 16

# Row 846
$ws.Range("A846").Value = 10
$ws.Range("B846").Value = "Vega Modelo de Temuco"
$ws.Range("C846").Value = "La Araucanía"
$ws.Range("D846").Value = 44214
$ws.Range("E846").Value = 9
$ws.Range("F846").Value = "Fruta"
$ws.Range("G846").Value = 100108
$ws.Range("H846").Value = "Tropicales y subtropicales"
$ws.Range("I846").Value = 100108005
$ws.Range("J846").Value = "Piña"
$ws.Range("K846").Value = "Caramelo"
$ws.Range("L846").Value = "Primera"
$ws.Range("M846").Value = 70
$ws.Range("N846").Value = 16000
$ws.Range("O846").Value = 17000
$ws.Range("P846").Value = 16571
$ws.Range("Q846").Value = "`$/caja 12 unidades"
$ws.Range("R846").Value = "Ecuador"
$ws.Range("S846").Value = 1381
$ws.Range("T846").Value = 12

# Row 847
$ws.Range("A847").Value = 10
$ws.Range("B847").Value = "Vega Modelo de Temuco"
$ws.Range("C847").Value = "La Araucanía"
$ws.Range("D847").Value = 44214
$ws.Range("E847").Value = 9
$ws.Range("F847").Value = "Fruta"
$ws.Range("G847").Value = 100108
$ws.Range("H847").Value = "Tropicales y subtropicales"
$ws.Range("I847").Value = 100108005
$ws.Range("J847").Value = "Piña"
$ws.Range("K847").Value = "Caramelo"
$ws.Range("L847").Value = "Segunda"
$ws.Range("M847").Value = 75
$ws.Range("N847").Value = 16000
$ws.Range("O847").Value = 16000
$ws.Range("P847").Value = 16000
$ws.Range("Q847").Value = "`$/caja 14 unidades"
$ws.Range("R847").Value = "Ecuador"
$ws.Range("S847").Value = 1143
$ws.Range("T847").Value = 14

# Row 848
$ws.Range("A848").Value = 10
$ws.Range("B848").Value = "Vega Modelo de Temuco"
$ws.Range("C848").Value = "La Araucanía"
$ws.Range("D848").Value = 44895
$ws.Range("E848").Value = 9
$ws.Range("F848").Value = "Fruta"
$ws.Range("G848").Value = 100108
$ws.Range("H848").Value = "Tropicales y subtropicales"
$ws.Range("I848").Value = 100108005
$ws.Range("J848").Value = "Piña"
$ws.Range("K848").Value = "Caramelo"
$ws.Range("L848").Value = "Segunda"
$ws.Range("M848").Value = 25
$ws.Range("N848").Value = 32000
$ws.Range("O848").Value = 32000
$ws.Range("P848").Value = 32000
$ws.Range("Q848").Value = "`$/caja 14 unidades"
$ws.Range("R848").Value = "Ecuador"
$ws.Range("S848").Value = 2286
$ws.Range("T848").Value = 14

# Row 849
$ws.Range("A849").Value = 10
$ws.Range("B849").Value = "Vega Modelo de Temuco"
$ws.Range("C849").Value = "La Araucanía"
$ws.Range("D849").Value = 44895
$ws.Range("E849").Value = 9
$ws.Range("F849").Value = "Fruta"
$ws.Range("G849").Value = 100108
$ws.Range("H849").Value = "Tropicales y subtropicales"
$ws.Range("I849").Value = 100108005
$ws.Range("J849").Value = "Piña"
$ws.Range("K849").Value = "Caramelo"
$ws.Range("L849").Value = "Segunda"
$ws.Range("M849").Value = 35
$ws.Range("N849").Value = 32000
$ws.Range("O849").Value = 32000
$ws.Range("P849").Value = 32000
$ws.Range("Q849").Value = "`$/caja 14 unidades"
$ws.Range("R849").Value = "Ecuador"
$ws.Range("S849").Value = 2286
$ws.Range("T849").Value = 14

# Row 850
$ws.Range("A850").Value = 10
$ws.Range("B850").Value = "Vega Modelo de Temuco"
$ws.Range("C850").Value = "La Araucanía"
$ws.Range("D850").Value = 45135
$ws.Range("E850").Value = 9
$ws.Range("F850").Value = "Fruta"
$ws.Range("G850").Value = 100108
$ws.Range("H850").Value = "Tropicales y subtropicales"
$ws.Range("I850").Value = 100108005
$ws.Range("J850").Value = "Piña"
$ws.Range("K850").Value = "Caramelo"
$ws.Range("L850").Value = "Primera"
$ws.Range("M850").Value = 130
$ws.Range("N850").Value = 25000
$ws.Range("O850").Value = 27000
$ws.Range("P850").Value = 26231
$ws.Range("Q850").Value = "`$/caja 12 unidades"
$ws.Range("R850").Value = "Ecuador"
$ws.Range("S850").Value = 2186
$ws.Range("T850").Value = 12

# Row 851
$ws.Range("A851").Value = 10
$ws.Range("B851").Value = "Vega Modelo de Temuco"
$ws.Range("C851").Value = "La Araucanía"
$ws.Range("D851").Value = 44567
$ws.Range("E851").Value = 9
$ws.Range("F851").Value = "Fruta"
$ws.Range("G851").Value = 100108
$ws.Range("H851").Value = "Tropicales y subtropicales"
$ws.Range("I851").Value = 100108005
$ws.Range("J851").Value = "Piña"
$ws.Range("K851").Value = "Caramelo"
$ws.Range("L851").Value = "Primera"
$ws.Range("M851").Value = 65
$ws.Range("N851").Value = 17000
$ws.Range("O851").Value = 18000
$ws.Range("P851").Value = 17385
$ws.Range("Q851").Value = "`$/caja 12 unidades"
$ws.Range("R851").Value = "Ecuador"
$ws.Range("S851").Value = 1449
$ws.Range("T851").Value = 12

# Row 852
$ws.Range("A852").Value = 10
$ws.Range("B852").Value = "Vega Modelo de Temuco"
$ws.Range("C852").Value = "La Araucanía"
$ws.Range("D852").Value = 44567
$ws.Range("E852").Value = 9
$ws.Range("F852").Value = "Fruta"
$ws.Range("G852").Value = 100108
$ws.Range("H852").Value = "Tropicales y subtropicales"
$ws.Range("I852").Value = 100108005
$ws.Range("J852").Value = "Piña"
$ws.Range("K852").Value = "Caramelo"
$ws.Range("L852").Value = "Segunda"
$ws.Range("M852").Value = 75
$ws.Range("N852").Value = 17000
$ws.Range("O852").Value = 17000
$ws.Range("P852").Value = 17000
$ws.Range("Q852").Value = "`$/caja 14 unidades"
$ws.Range("R852").Value = "Ecuador"
$ws.Range("S852").Value = 1214
$ws.Range("T852").Value = 14

# Row 853
$ws.Range("A853").Value = 10
$ws.Range("B853").Value = "Vega Modelo de Temuco"
$ws.Range("C853").Value = "La Araucanía"
$ws.Range("D853").Value = 45121
$ws.Range("E853").Value = 9
$ws.Range("F853").Value = "Fruta"
$ws.Range("G853").Value = 100108
$ws.Range("H853").Value = "Tropicales y subtropicales"
$ws.Range("I853").Value = 100108005
$ws.Range("J853").Value = "Piña"
$ws.Range("K853").Value = "Caramelo"
$ws.Range("L853").Value = "Primera"
$ws.Range("M853").Value = 80
$ws.Range("N853").Value = 25000
$ws.Range("O853").Value = 25000
$ws.Range("P853").Value = 25000
$ws.Range("Q853").Value = "`$/caja 12 unidades"
$ws.Range("R853").Value = "Ecuador"
$ws.Range("S853").Value = 2083
$ws.Range("T853").Value = 12

# Row 854
$ws.Range("A854").Value = 10
$ws.Range("B854").Value = "Vega Modelo de Temuco"
$ws.Range("C854").Value = "La Araucanía"
$ws.Range("D854").Value = 44203
$ws.Range("E854").Value = 9
$ws.Range("F854").Value = "Fruta"
$ws.Range("G854").Value = 100108
$ws.Range("H854").Value = "Tropicales y subtropicales"
$ws.Range("I854").Value = 100108005
$ws.Range("J854").Value = "Piña"
$ws.Range("K854").Value = "Caramelo"
$ws.Range("L854").Value = "Primera"
$ws.Range("M854").Value = 110
$ws.Range("N854").Value = 16000
$ws.Range("O854").Value = 17000
$ws.Range("P854").Value = 16591
$ws.Range("Q854").Value = "`$/caja 12 unidades"
$ws.Range("R854").Value = "Ecuador"
$ws.Range("S854").Value = 1383
$ws.Range("T854").Value = 12

# Row 855
$ws.Range("A855").Value = 10
$ws.Range("B855").Value = "Vega Modelo de Temuco"
$ws.Range("C855").Value = "La Araucanía"
$ws.Range("D855").Value = 44203
$ws.Range("E855").Value = 9
$ws.Range("F855").Value = "Fruta"
$ws.Range("G855").Value = 100108
$ws.Range("H855").Value = "Tropicales y subtropicales"
$ws.Range("I855").Value = 100108005
$ws.Range("J855").Value = "Piña"
$ws.Range("K855").Value = "Caramelo"
$ws.Range("L855").Value = "Segunda"
$ws.Range("M855").Value = 55
$ws.Range("N855").Value = 16000
$ws.Range("O855").Value = 16000
$ws.Range("P855").Value = 16000
$ws.Range("Q855").Value = "`$/caja 14 unidades"
$ws.Range("R855").Value = "Ecuador"
$ws.Range("S855").Value = 1143
$ws.Range("T855").Value = 14

# Row 856
$ws.Range("A856").Value = 10
$ws.Range("B856").Value = "Vega Modelo de Temuco"
$ws.Range("C856").Value = "La Araucanía"
$ws.Range("D856").Value = 44977
$ws.Range("E856").Value = 9
$ws.Range("F856").Value = "Fruta"
$ws.Range("G856").Value = 100108
$ws.Range("H856").Value = "Tropicales y subtropicales"
$ws.Range("I856").Value = 100108005
$ws.Range("J856").Value = "Piña"
$ws.Range("K856").Value = "Caramelo"
$ws.Range("L856").Value = "Primera"
$ws.Range("M856").Value = 150
$ws.Range("N856").Value = 26000
$ws.Range("O856").Value = 26000
$ws.Range("P856").Value = 26000
$ws.Range("Q856").Value = "`$/caja 12 unidades"
$ws.Range("R856").Value = "Ecuador"
$ws.Range("S856").Value = 2167
$ws.Range("T856").Value = 12

# Row 857
$ws.Range("A857").Value = 10
$ws.Range("B857").Value = "Vega Modelo de Temuco"
$ws.Range("C857").Value = "La Araucanía"
$ws.Range("D857").Value = 44977
$ws.Range("E857").Value = 9
$ws.Range("F857").Value = "Fruta"
$ws.Range("G857").Value = 100108
$ws.Range("H857").Value = "Tropicales y subtropicales"
$ws.Range("I857").Value = 100108005
$ws.Range("J857").Value = "Piña"
$ws.Range("K857").Value = "Caramelo"
$ws.Range("L857").Value = "Segunda"
$ws.Range("M857").Value = 100
$ws.Range("N857").Value = 26000
$ws.Range("O857").Value = 26000
$ws.Range("P857").Value = 26000
$ws.Range("Q857").Value = "`$/caja 14 unidades"
$ws.Range("R857").Value = "Ecuador"
$ws.Range("S857").Value = 1857
$ws.Range("T857").Value = 14

# Row 858
$ws.Range("A858").Value = 10
$ws.Range("B858").Value = "Vega Modelo de Temuco"
$ws.Range("C858").Value = "La Araucanía"
$ws.Range("D858").Value = 44818
$ws.Range("E858").Value = 9
$ws.Range("F858").Value = "Fruta"
$ws.Range("G858").Value = 100108
$ws.Range("H858").Value = "Tropicales y subtropicales"
$ws.Range("I858").Value = 100108005
$ws.Range("J858").Value = "Piña"
$ws.Range("K858").Value = "Caramelo"
$ws.Range("L858").Value = "Primera"
$ws.Range("M858").Value = 160
$ws.Range("N858").Value = 24000
$ws.Range("O858").Value = 25000
$ws.Range("P858").Value = 24500
$ws.Range("Q858").Value = "`$/caja 12 unidades"
$ws.Range("R858").Value = "Ecuador"
$ws.Range("S858").Value = 2042
$ws.Range("T858").Value = 12

# Row 859
$ws.Range("A859").Value = 10
$ws.Range("B859").Value = "Vega Modelo de Temuco"
$ws.Range("C859").Value = "La Araucanía"
$ws.Range("D859").Value = 44293
$ws.Range("E859").Value = 9
$ws.Range("F859").Value = "Fruta"
$ws.Range("G859").Value = 100108
$ws.Range("H859").Value = "Tropicales y subtropicales"
$ws.Range("I859").Value = 100108005
$ws.Range("J859").Value = "Piña"
$ws.Range("K859").Value = "Caramelo"
$ws.Range("L859").Value = "Primera"
$ws.Range("M859").Value = 100
$ws.Range("N859").Value = 16000
$ws.Range("O859").Value = 17000
$ws.Range("P859").Value = 16500
$ws.Range("Q859").Value = "`$/caja 12 unidades"
$ws.Range("R859").Value = "Ecuador"
$ws.Range("S859").Value = 1375
$ws.Range("T859").Value = 12

# Row 860
$ws.Range("A860").Value = 10
$ws.Range("B860").Value = "Vega Modelo de Temuco"
$ws.Range("C860").Value = "La Araucanía"
$ws.Range("D860").Value = 44293
$ws.Range("E860").Value = 9
$ws.Range("F860").Value = "Fruta"
$ws.Range("G860").Value = 100108
$ws.Range("H860").Value = "Tropicales y subtropicales"
$ws.Range("I860").Value = 100108005
$ws.Range("J860").Value = "Piña"
$ws.Range("K860").Value = "Caramelo"
$ws.Range("L860").Value = "Segunda"
$ws.Range("M860").Value = 50
$ws.Range("N860").Value = 17000
$ws.Range("O860").Value = 17000
$ws.Range("P860").Value = 17000
$ws.Range("Q860").Value = "`$/caja 14 unidades"
$ws.Range("R860").Value = "Ecuador"
$ws.Range("S860").Value = 1214
$ws.Range("T860").Value = 14

# Row 861
$ws.Range("A861").Value = 10
$ws.Range("B861").Value = "Vega Modelo de Temuco"
$ws.Range("C861").Value = "La Araucanía"
$ws.Range("D861").Value = 44189
$ws.Range("E861").Value = 9
$ws.Range("F861").Value = "Fruta"
$ws.Range("G861").Value = 100108
$ws.Range("H861").Value = "Tropicales y subtropicales"
$ws.Range("I861").Value = 100108005
$ws.Range("J861").Value = "Piña"
$ws.Range("K861").Value = "Caramelo"
$ws.Range("L861").Value = "Primera"
$ws.Range("M861").Value = 130
$ws.Range("N861").Value = 17000
$ws.Range("O861").Value = 18000
$ws.Range("P861").Value = 17385
$ws.Range("Q861").Value = "`$/caja 12 unidades"
$ws.Range("R861").Value = "Ecuador"
$ws.Range("S861").Value = 1449
$ws.Range("T861").Value = 12
$ws.Range("D861").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 862
$ws.Range("A862").Value = 10
$ws.Range("B862").Value = "Vega Modelo de Temuco"
$ws.Range("C862").Value = "La Araucanía"
$ws.Range("D862").Value = 45100
$ws.Range("E862").Value = 9
$ws.Range("F862").Value = "Fruta"
$ws.Range("G862").Value = 100108
$ws.Range("H862").Value = "Tropicales y subtropicales"
$ws.Range("I862").Value = 100108005
$ws.Range("J862").Value = "Piña"
$ws.Range("K862").Value = "Caramelo"
$ws.Range("L862").Value = "Primera"
$ws.Range("M862").Value = 85
$ws.Range("N862").Value = 33000
$ws.Range("O862").Value = 33000
$ws.Range("P862").Value = 33000
$ws.Range("Q862").Value = "`$/caja 12 unidades"
$ws.Range("R862").Value = "Ecuador"
$ws.Range("S862").Value = 2750
$ws.Range("T862").Value = 12
$ws.Range("D862").NumberFormat = "YYYY-MM-DD HH:MM:SS"
